# Update "list articles" and "pipeline" sheets with corrected/normalized
# article titles and new review statuses, per commit "update files from 4 to 9 + metadata"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("list articles")
$ws2 = $wb.Worksheets.Item("pipeline")

# --- Sheet "list articles" ---

# Row 1 (header) - status column gets a value
$ws1.Range("C1").Value = "manual rewriting"

# Row 2
$ws1.Range("B2").Value = "TUNDRA: A multilingual corpus of found data for TTS research created"
$ws1.Range("C2").Value = "ok"

# Row 3
$ws1.Range("B3").Value = "Current trends in multilingual speech processing"
$ws1.Range("C3").Value = "don't need"

# Row 4
$ws1.Range("B4").Value = "Learning to speak fluently in a foreign language:"
$ws1.Range("C4").Value = "ok"

# Row 5
$ws1.Range("B5").Value = "One model, many languages: meta-learning for multilingual Text-to-Speech"
$ws1.Range("C5").Value = "ok"

# Row 6
$ws1.Range("B6").Value = "Non-autoregressive neural Text-to-Speech"
$ws1.Range("C6").Value = "ok"

# Row 7
$ws1.Range("B7").Value = "Directly modeling speech waveforms by neural networks"
$ws1.Range("C7").Value = "ok"

# Row 8 - title cleared
$ws1.Range("B8").ClearContents()

# Row 9
$ws1.Range("B9").Value = "Multilingual text Aanalysis for Text-to-Speech synthesis"
$ws1.Range("C9").Value = "ok"

# Row 10
$ws1.Range("B10").Value = "Text-to-Speech conveersion with neurzl networks"
$ws1.Range("C10").Value = "ok"

# Row 11
$ws1.Range("B11").Value = "MERLIN: An open source neural network speech synthesis system"

# Row 12
$ws1.Range("B12").Value = "Grapheme-to-Phoneme conversion with convoltional neural networks"

# Row 13, 14 - unchanged (still blank)

# Row 15
$ws1.Range("B15").Value = "GlobalPhone: A multilingual Text & Speech database in 20 Languages"

# Row 16
$ws1.Range("B16").Value = "Learning to speak fluently in a foreign language:"

# Row 17
$ws1.Range("B17").Value = " "

# Row 18
$ws1.Range("B18").Value = "Multilingual Text-to-Speech synthesis"

# Row 19
$ws1.Range("B19").Value = "MELLOTRON: Multispeaker expressive voice synthesis by conditioning"

# Row 20
$ws1.Range("B20").Value = "Text Preprocessing for Speech Synthesis"

# Row 21
$ws1.Range("B21").Value = "Unsupervised and lightly-supervised learning for rapid construction of TTS"

# Row 22
$ws1.Range("B22").Value = "Statistical parametric speech synthesis"

# Row 23
$ws1.Range("B23").Value = "TACOTRON: Toward end-to-end speech synthesis"
$ws1.Range("C23").Value = "don't need"

# Update selected cell on sheet "list articles"
$ws1.Activate()
$ws1.Range("B11").Select() | Out-Null

# --- Sheet "pipeline" (values unchanged, only shared-string reindex upstream) ---
$ws2.Range("A1").Value = "steps"
$ws2.Range("B2").Value = "create text files (corpus)"
$ws2.Range("B3").Value = "build lexicon from corpus (using Termsuite -terminology extraction- and manually)"
$ws2.Range("B4").Value = "build rules based system (extract terms from lexicon)"
$ws2.Range("B5").Value = "run tagger (with IOB tagset)"

$wb.Save()
